$d = $word.ActiveDocument

# Locate the paragraph that holds "LOB1036: Geometria Analítica (Requisito fraco)".
# The four paragraphs that follow it are the ones being removed by this edit:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) an empty paragraph
#   4) an empty paragraph with PageBreakBefore + left alignment
# They are replaced by nothing, so the paragraph right after LOB1036 becomes
# the empty "pStyle only" paragraph that used to sit two spots further down
# (the one that still precedes the remaining PageBreakBefore paragraph).

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1036*Geometria Anal*") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find the LOB1036 paragraph"
}

$first = $anchor.Next()
$last = $first
for ($k = 1; $k -lt 4; $k++) {
    $last = $last.Next()
}

$rangeToDelete = $d.Range($first.Range.Start, $last.Range.End)
$rangeToDelete.Delete()
